# Update two links in the "guides" sheet of the data-mgmt resources workbook:
#  1. The "DIME Analytics" row's link text/target changes from the old
#     "Dimewiki" (Primary_Data_Collection) page to the new "DIME Wiki"
#     (Main_Page) page.
#  2. A brand-new "Responsible Data" guide entry (linking to "The Handbook of
#     the Modern Development Specialist") is added, in its correct
#     alphabetically-sorted position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("guides")

# Remove all existing hyperlink bookkeeping up front - we will recreate every
# hyperlink at the end once all the cells/rows are in their final place, so
# that each hyperlink ref points at the right (possibly shifted) cell.
$ws.Cells.Hyperlinks.Delete()

# --- 1. Update the DIME Analytics row (row 3) -----------------------------
# A3 (the org name "DIME Analytics") is unchanged; only the link in B3 moves
# from the old Dimewiki page to the new DIME Wiki page.
$ws.Range("B3").Value = '<a href="https://dimewiki.worldbank.org/Main_Page">DIME Wiki</a>'

# --- 2. Insert the new "Responsible Data" row ------------------------------
# It belongs alphabetically between "Lewis, C." (row 9) and
# "Reynolds, T., Schatschneider, C. & Logan, J." (old row 10), i.e. the new
# row 10.
$ws.Rows.Item(10).Insert()
$ws.Range("A10").Value = "Responsible Data"
$ws.Range("B10").Value = '<a href="https://responsibledata.io/resources/handbook/">The Handbook of the Modern Development Specialist</a>'

# --- 3. Keep the sortState / sortCondition ranges in sync with the new size
$dataRows = $ws.Cells.Item(1, 1).CurrentRegion.Rows.Count
$lastRow = $dataRows
$ws.Range("A2:A$lastRow").Sort($ws.Range("A2:A$lastRow"))

# --- 4. Re-create every hyperlink, in the same order as the original file,
#     so relationship ids come out as rId1..rId7 and every ref points at the
#     (possibly shifted) correct cell.
$ws.Hyperlinks.Add($ws.Range("B4"), "https://rdmkit.elixir-europe.org/", [System.Type]::Missing, [System.Type]::Missing, "https://rdmkit.elixir-europe.org/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B12"), "https://the-turing-way.netlify.app/welcome", [System.Type]::Missing, [System.Type]::Missing, "https://the-turing-way.netlify.app/welcome") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B11"), "https://the-turing-way.netlify.app/welcome", [System.Type]::Missing, [System.Type]::Missing, "https://the-turing-way.netlify.app/welcome") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B8"), "https://the-turing-way.netlify.app/welcome", [System.Type]::Missing, [System.Type]::Missing, "https://the-turing-way.netlify.app/welcome") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "https://the-turing-way.netlify.app/welcome", [System.Type]::Missing, [System.Type]::Missing, "https://the-turing-way.netlify.app/welcome") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B5"), "https://rdmkit.elixir-europe.org/", [System.Type]::Missing, [System.Type]::Missing, "https://rdmkit.elixir-europe.org/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B10"), "https://rdmkit.elixir-europe.org/", [System.Type]::Missing, [System.Type]::Missing, "https://rdmkit.elixir-europe.org/") | Out-Null

# Re-apply the display text for every hyperlinked cell (adding the hyperlink
# with a TextToDisplay above forces the cell text to the url; put back the
# real cell contents now) and strip the auto-applied "Hyperlink" style so the
# cells keep the workbook's normal (unstyled) look.
$ws.Range("B4").Value = '<a href="https://rdmkit.elixir-europe.org/">Research Data Management Kit</a>'
$ws.Range("B12").Value = '<a href="https://the-turing-way.netlify.app/welcome">Handbook</a>'
$ws.Range("B11").Value = '<a href="https://figshare.com/articles/preprint/The_Basics_of_Data_Management/13215350">The Basics of Data Management</a>'
$ws.Range("B8").Value = '<a href="https://www.povertyactionlab.org/resource/introduction-randomized-evaluations">Research Resources</a>'
$ws.Range("B3").Value = '<a href="https://dimewiki.worldbank.org/Main_Page">DIME Wiki</a>'
$ws.Range("B5").Value = '<a href="https://www.sjsu.edu/research/docs/irb-data-management-handbook.pdf">Data Management Handbook for Human Subjects Research</a>'
$ws.Range("B10").Value = '<a href="https://responsibledata.io/resources/handbook/">The Handbook of the Modern Development Specialist</a>'

$ws.Range("B3:B12").Style = "Normal"

# --- 5. Restore the sheet/window selection state ---------------------------
$ws.Activate()
$ws.Range("C2").Select() | Out-Null
